# Auto update: 2025-11-29 03:48:30
# Updates the 미장_양자_분석 (US quantum computing stock analysis) sheet with the
# latest day's data. Row 2 now holds D-Wave (QBTS) and row 3 now holds
# International Business Machines (IBM); their numeric metrics are refreshed,
# and rows 4-5 (Rigetti/RGTI, IonQ/IONQ) get refreshed score/macro values too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D-Wave Quantum Inc. / QBTS
$ws.Cells.Item(2, 2).Value = "D-Wave Quantum Inc."
$ws.Cells.Item(2, 3).Value = "QBTS"
$ws.Range("D2").Value = 22.67
$ws.Range("F2").Value = 10.53
$ws.Range("G2").Value = 20
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 62.1
$ws.Range("N2").Value = 85.82376350509293

# Row 3: International Business Machines / IBM
$ws.Cells.Item(3, 2).Value = "International Business Machines"
$ws.Cells.Item(3, 3).Value = "IBM"
$ws.Range("D3").Value = 308.58
$ws.Range("F3").Value = 6.26
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 61.7
$ws.Range("N3").Value = 85.82376350509293

# Row 4: Rigetti Computing, Inc. / RGTI - score + macro refresh only
$ws.Range("K4").Value = 59.7
$ws.Range("N4").Value = 85.82376350509293

# Row 5: IonQ, Inc. / IONQ
$ws.Range("D5").Value = 49.3
$ws.Range("F5").Value = 20.24
$ws.Range("K5").Value = 57.1
$ws.Range("N5").Value = 85.82376350509293
